$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 120: code=10113 / postal code (ara)
$ws.Cells.Item(120, 1).Value = 10113
$ws.Cells.Item(120, 2).Value = 10113
$ws.Cells.Item(120, 3).Value = 5
$ws.Cells.Item(120, 4).Value = "الرمز البريدي"
$ws.Cells.Item(120, 5).Value = "BNMR"
$ws.Cells.Item(120, 6).Value = "ara"
$ws.Cells.Item(120, 7).Value = $true
$ws.Cells.Item(120, 8).Value = "superadmin"
$ws.Cells.Item(120, 9).Value = "now()"

# Row 121: code=10114 / postal code (ara)
$ws.Cells.Item(121, 1).Value = 10114
$ws.Cells.Item(121, 2).Value = 10114
$ws.Cells.Item(121, 3).Value = 5
$ws.Cells.Item(121, 4).Value = "الرمز البريدي"
$ws.Cells.Item(121, 5).Value = "BNMR"
$ws.Cells.Item(121, 6).Value = "ara"
$ws.Cells.Item(121, 7).Value = $true
$ws.Cells.Item(121, 8).Value = "superadmin"
$ws.Cells.Item(121, 9).Value = "now()"

# Update the selection to match the author's saved view (rows below the new data)
[void]$ws.Rows("122:1048576").Select()
